# Etienne's update with substantial changes...
# Applies the KES-MPS.xlsx edit: new "KES Det. unit FOV / MPS" row, a new
# mu / Se geometry block to the right of the Brent sheet, a couple of new
# data points (I13/J13), plus the cosmetic zoom / selection changes that
# came along with the edit.

$wb = $excel.ActiveWorkbook

$brent = $wb.Worksheets.Item("Brent")
$lin   = $wb.Worksheets.Item("Lin")

# ---------------------------------------------------------------------
# Brent sheet
# ---------------------------------------------------------------------

# New data points next to the existing "Effective thickness" row.
$brent.Range("I13").Value = 17.5
$brent.Range("J13").Value = 15.6

# New row 17: "KES Det. unit FOV / MPS" ratios (mirrors row 16's pattern).
$brent.Range("A17").Value = "KES Det. unit FOV / MPS"
$brent.Range("B17").Formula = "=C13/B13"
$brent.Range("I17").Formula = "=J13/I13"

# New geometry mini-table to the right (U:X), next to the footnote rows.
$brent.Range("U18").Value = "µ (mm-1)"
$brent.Range("V18").Formula = "=4.4*10^(-2)*18/10"
$brent.Range("X18").Formula = "=EXP(-V18*5)"

$brent.Range("U19").Value = "Se (mm)"
$brent.Range("V19").Formula = "=LN(2)/(V18*TAN(60*3.14159/180))"

# Cosmetic: zoom level dropped from 110% to 90%, and the selection moved.
$brent.Activate()
$excel.ActiveWindow.Zoom = 90
$brent.Range("E1").Select()

# ---------------------------------------------------------------------
# Lin sheet
# ---------------------------------------------------------------------

$lin.Activate()
$excel.ActiveWindow.Zoom = 90
$lin.Range("D8").Select()
